# Insert a new data row just before the existing row 427. This shifts the
# current rows 427-525 down to 428-526 (preserving all their data/formatting)
# and leaves a fresh blank row 427 (inheriting the row's formatting, notably
# the date style on column D) which we then populate with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("427:427").Insert()

$ws.Range("A427").Value = 3
$ws.Range("B427").Value = "Femacal de La Calera"
$ws.Range("C427").Value = "Coquimbo"
$ws.Range("D427").Value = 44798
$ws.Range("E427").Value = 5
$ws.Range("F427").Value = 100112021
$ws.Range("G427").Value = "Ají"
$ws.Range("H427").Value = "Inferno"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 78
$ws.Range("K427").Value = 11000
$ws.Range("L427").Value = 12000
$ws.Range("M427").Value = 11513
$ws.Range("N427").Value = "$/caja 10 kilos"
$ws.Range("O427").Value = "Región de Arica y Parinacota"
$ws.Range("P427").Value = 1151
$ws.Range("Q427").Value = 10
$ws.Range("R427").Value = "Hortaliza"
